$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Unit Price" column (I). This shifts
#    Preservation Method (J->I) and Due Date (K->J) left.
$ws.Range("I1").EntireColumn.Delete()

# 2) Correct two preservation-method codes that were missing their
#    "M" prefix (visible now that they live in column I).
$ws.Range('I9').Value = 'M41'
$ws.Range('I15').Value = 'M10'

# 3) Append the newly awarded contracts (rows 16-21). Force the
#    cells to text first so values like dates, NSNs, and currency
#    strings are stored verbatim instead of being auto-converted.
$row16 = $ws.Range("A16:J16")
$row16.NumberFormat = "@"
$ws.Range('A16').Value = '3/25/2019'
$ws.Range('B16').Value = 'SPE7M2-19-V-1046'
$ws.Range('C16').Value = '115'
$ws.Range('D16').Value = '$32,557.65'
$ws.Range('E16').Value = '5330010778314'
$ws.Range('F16').Value = 'PACKING ASSORTMENT,PREFORMED'
$ws.Range('G16').Value = 'Munters'
$ws.Range('H16').Value = '43125-08'
$ws.Range('I16').Value = 'M33'
$ws.Range('J16').Value = '2019 SEP 03'
$row16.ClearFormats()

$row17 = $ws.Range("A17:J17")
$row17.NumberFormat = "@"
$ws.Range('A17').Value = '3/25/2019'
$ws.Range('B17').Value = 'SPE7L1-19-V-4818'
$ws.Range('C17').Value = '2'
$ws.Range('D17').Value = '$2,224.00'
$ws.Range('E17').Value = '3020016570251'
$ws.Range('F17').Value = 'GEAR,ANTIBACKLASH,H'
$ws.Range('G17').Value = 'KTSDI'
$ws.Range('H17').Value = '510251-163'
$ws.Range('I17').Value = 'M41'
$ws.Range('J17').Value = '2019 JUN 28'
$row17.ClearFormats()

$row18 = $ws.Range("A18:J18")
$row18.NumberFormat = "@"
$ws.Range('A18').Value = '3/25/2019'
$ws.Range('B18').Value = 'SPE7M0-19-V-6074'
$ws.Range('C18').Value = '6'
$ws.Range('D18').Value = '$12,650.22'
$ws.Range('E18').Value = '5930013674492'
$ws.Range('F18').Value = 'SWITCH,FLOW'
$ws.Range('G18').Value = 'GEMS'
$ws.Range('H18').Value = '139644'
$ws.Range('I18').Value = 'CP'
$ws.Range('J18').Value = '2019 SEP 03'
$row18.ClearFormats()

$row19 = $ws.Range("A19:J19")
$row19.NumberFormat = "@"
$ws.Range('A19').Value = '3/25/2019'
$ws.Range('B19').Value = 'SPE7M0-19-V-6059'
$ws.Range('C19').Value = '6'
$ws.Range('D19').Value = '$6,219.84'
$ws.Range('E19').Value = '5930014842487'
$ws.Range('F19').Value = 'SWITCH,FLOW'
$ws.Range('G19').Value = 'GEMS'
$ws.Range('H19').Value = '159297'
$ws.Range('I19').Value = 'CP'
$ws.Range('J19').Value = '2019 SEP 03'
$row19.ClearFormats()

$row20 = $ws.Range("A20:J20")
$row20.NumberFormat = "@"
$ws.Range('A20').Value = '3/25/2019'
$ws.Range('B20').Value = 'SPE7MC-19-V-5914'
$ws.Range('C20').Value = '13'
$ws.Range('D20').Value = '$4,654.00 '
$ws.Range('E20').Value = '5935014967270'
$ws.Range('F20').Value = 'BACKSHELL,ELECTRICAL CONNECTOR'
$ws.Range('G20').Value = 'Glenair'
$ws.Range('H20').Value = '447HT325XW2519'
$ws.Range('I20').Value = 'CP'
$ws.Range('J20').Value = '2019 SEP 03'
$row20.ClearFormats()

$row21 = $ws.Range("A21:J21")
$row21.NumberFormat = "@"
$ws.Range('A21').Value = '3/25/2019'
$ws.Range('B21').Value = 'SPE7M1-19-V-5446'
$ws.Range('C21').Value = '2'
$ws.Range('D21').Value = '$9,545.28 '
$ws.Range('E21').Value = '6110014091404'
$ws.Range('F21').Value = 'CONTROLLER,MOTOR'
$ws.Range('G21').Value = 'Morpac'
$ws.Range('H21').Value = '70000-4'
$ws.Range('I21').Value = 'CP'
$ws.Range('J21').Value = '2019 SEP 03'
$row21.ClearFormats()

# 4) Restore the active selection to match the edited area.
$ws.Range("I10").Select()
